# ---------------------------------------------------------------------------
# Edit: 
#   1. Slide 6 table ("SOURCES OF FINANCE") switches its table style from
#      the custom "Table_0" style ({3D2EC5C2-88CF-4051-B577-94BBD11D5F38})
#      to the built-in style {F9392547-D375-465A-AC56-47E2994FE49A}.
#   2. The presentation's main theme (ppt/theme/theme1.xml, used by the
#      slide master) has its colour palette swapped from the custom
#      "Integral" palette to the standard Office palette (the font scheme
#      and format scheme were already identical between the two themes
#      present in the deck, so only the 12 theme colours actually change).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ---------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F9392547-D375-465A-AC56-47E2994FE49A}")
    }
}

# --- 2. Theme colours (Integral -> Office) on the slide master's theme ----
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
